$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 32
$ws.Range("H32").Value = 2754
$ws.Range("I32").Value = 2754
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2754
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2428
$ws.Range("N32").ClearContents()
# Row 33
$ws.Range("H33").Value = 2129.25
$ws.Range("I33").Value = 1658.85
$ws.Range("J33").Value = 3305.25
$ws.Range("K33").Value = 1658.85
$ws.Range("L33").Value = 3305.25
$ws.Range("M33").Value = -1429.85
$ws.Range("N33").Value = -3763.25
# Row 98
$ws.Range("H98").Value = 50002284
$ws.Range("I98").Value = 52633776
$ws.Range("K98").Value = 52633776
$ws.Range("M98").Value = -52632278
# Row 116
$ws.Range("H116").Value = 6302.647
$ws.Range("I116").Value = 5917.909
$ws.Range("J116").Value = 7008
$ws.Range("K116").Value = 5917.909
$ws.Range("L116").Value = 7008
$ws.Range("M116").Value = -2475.909
$ws.Range("N116").Value = -13892
# Row 122
$ws.Range("H122").Value = 50002284
$ws.Range("I122").Value = 52633776
$ws.Range("K122").Value = 157901328
$ws.Range("M122").Value = -157898878
# Row 131
$ws.Range("H131").Value = 3711.5293
$ws.Range("I131").Value = 1645.091
$ws.Range("K131").Value = 4935.272999999999
$ws.Range("M131").Value = 104.7270000000008
# Row 132
$ws.Range("H132").Value = 6048.25
$ws.Range("I132").Value = 6094.5
$ws.Range("J132").Value = 5847.8335
$ws.Range("K132").Value = 18283.5
$ws.Range("L132").Value = 17543.5005
$ws.Range("M132").Value = -15753.5
$ws.Range("N132").Value = -22603.5005
# Row 137
$ws.Range("H137").Value = 3511.4075
$ws.Range("I137").Value = 2512.36
$ws.Range("J137").Value = 15999.5
$ws.Range("K137").Value = 7537.08
$ws.Range("L137").Value = 47998.5
$ws.Range("M137").Value = -4987.08
$ws.Range("N137").Value = -53098.5
# Row 138
$ws.Range("H138").Value = 2475.8271
$ws.Range("I138").Value = 1503.95
$ws.Range("J138").Value = 2794.4753
$ws.Range("K138").Value = 4511.85
$ws.Range("L138").Value = 8383.4259
$ws.Range("M138").Value = 628.1499999999996
$ws.Range("N138").Value = -18663.4259
# Row 141
$ws.Range("H141").Value = 5396.2383
$ws.Range("I141").Value = 5586.05
$ws.Range("K141").Value = 16758.15
$ws.Range("M141").Value = -11578.15

$ws = $wb.Worksheets("ARM")
# Row 32
$ws.Range("H32").Value = 7056815.5
$ws.Range("I32").Value = 7826720.5
$ws.Range("J32").Value = 17687.285
$ws.Range("K32").Value = 7826720.5
$ws.Range("L32").Value = 17687.285
$ws.Range("M32").Value = -7826433.5
$ws.Range("N32").Value = -18261.285
# Row 61
$ws.Range("H61").Value = 8827285
$ws.Range("I61").Value = 6413375.5
$ws.Range("J61").Value = 35725136
$ws.Range("K61").Value = 6413375.5
$ws.Range("L61").Value = 35725136
$ws.Range("M61").Value = -6413163.5
$ws.Range("N61").Value = -35725560
# Row 74
$ws.Range("H74").Value = 18580554
$ws.Range("I74").Value = 41675652
$ws.Range("J74").Value = 1259228.1
$ws.Range("K74").Value = 41675652
$ws.Range("L74").Value = 1259228.1
$ws.Range("M74").Value = -41674778
$ws.Range("N74").Value = -1260976.1
# Row 77
$ws.Range("H77").Value = 18580554
$ws.Range("I77").Value = 41675652
$ws.Range("J77").Value = 1259228.1
$ws.Range("K77").Value = 208378260
$ws.Range("L77").Value = 6296140.5
$ws.Range("M77").Value = -208373892
$ws.Range("N77").Value = -6304876.5
# Row 117
$ws.Range("H117").Value = 24000
$ws.Range("J117").Value = 24000
$ws.Range("L117").Value = 24000
$ws.Range("N117").Value = -33178
# Row 119
$ws.Range("H119").Value = 524742.25
$ws.Range("J119").Value = 32989.668
$ws.Range("L119").Value = 32989.668
$ws.Range("N119").Value = -42665.668
# Row 121
$ws.Range("H121").Value = 110000
$ws.Range("J121").Value = 110000
$ws.Range("L121").Value = 110000
$ws.Range("N121").Value = -113494
# Row 132
$ws.Range("H132").Value = 2659.2727
$ws.Range("I132").Value = 1301.8096
$ws.Range("K132").Value = 3905.4288
$ws.Range("M132").Value = -1375.4288
# Row 136
$ws.Range("H136").Value = 8827285
$ws.Range("I136").Value = 6413375.5
$ws.Range("J136").Value = 35725136
$ws.Range("K136").Value = 19240126.5
$ws.Range("L136").Value = 107175408
$ws.Range("M136").Value = -19237576.5
$ws.Range("N136").Value = -107180508

$ws = $wb.Worksheets("BSM")
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
# Row 58
$ws.Range("H58").Value = 9000
$ws.Range("J58").Value = 9000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9588
# Row 134
$ws.Range("H134").Value = 334450.28
$ws.Range("I134").Value = 1005.96
$ws.Range("J134").Value = 2001671.8
$ws.Range("K134").Value = 3017.88
$ws.Range("L134").Value = 6005015.4
$ws.Range("M134").Value = -482.8800000000001
$ws.Range("N134").Value = -6010085.4

$ws = $wb.Worksheets("CRP")
# Row 31
$ws.Range("H31").Value = 941924.0600000001
$ws.Range("I31").Value = 9535.842000000001
$ws.Range("K31").Value = 9535.842000000001
$ws.Range("M31").Value = -9240.842000000001
# Row 34
$ws.Range("H34").Value = 941924.0600000001
$ws.Range("I34").Value = 9535.842000000001
$ws.Range("K34").Value = 9535.842000000001
$ws.Range("M34").Value = -9333.842000000001
# Row 58
$ws.Range("H58").Value = 2294.7026
$ws.Range("I58").Value = 1120.4138
$ws.Range("J58").Value = 6551.5
$ws.Range("K58").Value = 1120.4138
$ws.Range("L58").Value = 6551.5
$ws.Range("M58").Value = -917.4138
$ws.Range("N58").Value = -6957.5
# Row 132
$ws.Range("H132").Value = 1655.3673
$ws.Range("I132").Value = 1701.1428
$ws.Range("J132").Value = 1380.7142
$ws.Range("K132").Value = 5103.428400000001
$ws.Range("L132").Value = 4142.142599999999
$ws.Range("M132").Value = -2573.428400000001
$ws.Range("N132").Value = -9202.142599999999
# Row 134
$ws.Range("H134").Value = 2045.8837
$ws.Range("I134").Value = 1028.0857
$ws.Range("K134").Value = 3084.2571
$ws.Range("M134").Value = -549.2571000000003
# Row 136
$ws.Range("H136").Value = 2294.7026
$ws.Range("I136").Value = 1120.4138
$ws.Range("J136").Value = 6551.5
$ws.Range("K136").Value = 3361.2414
$ws.Range("L136").Value = 19654.5
$ws.Range("M136").Value = -811.2413999999999
$ws.Range("N136").Value = -24754.5

$ws = $wb.Worksheets("CUL")
# Row 42
$ws.Range("H42").Value = 7800
$ws.Range("I42").Value = 7800
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 23400
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -22866
$ws.Range("N42").ClearContents()
# Row 44
$ws.Range("H44").Value = 216
$ws.Range("I44").Value = 195
$ws.Range("J44").Value = 300
$ws.Range("K44").Value = 585
$ws.Range("L44").Value = 900
$ws.Range("M44").Value = -187
$ws.Range("N44").Value = -1696
# Row 50
$ws.Range("H50").Value = 421.19232
$ws.Range("I50").Value = 185
$ws.Range("J50").Value = 452
$ws.Range("K50").Value = 555
$ws.Range("L50").Value = 1356
$ws.Range("M50").Value = -74
$ws.Range("N50").Value = -2318
# Row 53
$ws.Range("H53").Value = 421.19232
$ws.Range("I53").Value = 185
$ws.Range("J53").Value = 452
$ws.Range("K53").Value = 555
$ws.Range("L53").Value = 1356
$ws.Range("M53").Value = -74
$ws.Range("N53").Value = -2318
# Row 113
$ws.Range("H113").Value = 1061.7273
$ws.Range("J113").Value = 1337.3334
$ws.Range("L113").Value = 4012.0002
$ws.Range("N113").Value = -8352.0002

$ws = $wb.Worksheets("GSM")
# Row 28
$ws.Range("H28").Value = 1500
$ws.Range("J28").Value = 1500
$ws.Range("L28").Value = 1500
$ws.Range("N28").Value = -1884
# Row 54
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -10780
# Row 107
$ws.Range("H107").Value = 446.69232
$ws.Range("I107").Value = 339.125
$ws.Range("K107").Value = 339.125
$ws.Range("M107").Value = 1580.875
# Row 132
$ws.Range("H132").Value = 13702185
$ws.Range("I132").Value = 17545888
$ws.Range("K132").Value = 52637664
$ws.Range("M132").Value = -52635134

$ws = $wb.Worksheets("LTW")
# Row 7
$ws.Range("H7").Value = 47169.707
$ws.Range("I7").Value = 4662.353
$ws.Range("J7").Value = 150401.86
$ws.Range("K7").Value = 4662.353
$ws.Range("L7").Value = 150401.86
$ws.Range("M7").Value = -4550.353
$ws.Range("N7").Value = -150625.86
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
# Row 22
$ws.Range("H22").Value = 1375.9375
$ws.Range("I22").Value = 1270.2
$ws.Range("J22").Value = 1424
$ws.Range("K22").Value = 1270.2
$ws.Range("L22").Value = 1424
$ws.Range("M22").Value = -975.2
$ws.Range("N22").Value = -2014
# Row 27
$ws.Range("H27").Value = 1375.9375
$ws.Range("I27").Value = 1270.2
$ws.Range("J27").Value = 1424
$ws.Range("K27").Value = 1270.2
$ws.Range("L27").Value = 1424
$ws.Range("M27").Value = -1163.2
$ws.Range("N27").Value = -1638
# Row 68
$ws.Range("H68").Value = 4812.4287
$ws.Range("J68").Value = 5128.7144
$ws.Range("L68").Value = 5128.7144
$ws.Range("N68").Value = -6626.7144
# Row 71
$ws.Range("H71").Value = 4812.4287
$ws.Range("J71").Value = 5128.7144
$ws.Range("L71").Value = 25643.572
$ws.Range("N71").Value = -33131.572
# Row 126
$ws.Range("H126").Value = 47169.707
$ws.Range("I126").Value = 4662.353
$ws.Range("J126").Value = 150401.86
$ws.Range("K126").Value = 13987.059
$ws.Range("L126").Value = 451205.58
$ws.Range("M126").Value = -11517.059
$ws.Range("N126").Value = -456145.58
# Row 132
$ws.Range("H132").Value = 843216.8
$ws.Range("I132").Value = 10781.546
$ws.Range("K132").Value = 32344.638
$ws.Range("M132").Value = -29814.638

$ws = $wb.Worksheets("WVR")
# Row 119
$ws.Range("H119").Value = 29989.25
$ws.Range("J119").Value = 29989.25
$ws.Range("L119").Value = 29989.25
$ws.Range("N119").Value = -39665.25
# Row 136
$ws.Range("H136").Value = 3010
$ws.Range("I136").Value = 2496.5334
$ws.Range("J136").Value = 5210.5713
$ws.Range("K136").Value = 7489.600199999999
$ws.Range("L136").Value = 15631.7139
$ws.Range("M136").Value = -4939.600199999999
$ws.Range("N136").Value = -20731.7139
